$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 17 ("linear programming" solvers slide): TextBox 3 (shape id 4)
# ---------------------------------------------------------------------
$s17 = $p.Slides.Item(17)
$box = $s17.Shapes.Item(4)

# Resize / reposition the text box (off x: 3956613 -> 2425701 EMU,
# ext cx: 8235387 -> 9766300 EMU; top/height unchanged).
$box.Left = 191.0
$box.Width = 769.0

$tr = $box.TextFrame.TextRange

# Merge "Some other " into the following run so the paragraph reads
# "Some other open-source optimization libraries " as a single run.
$full = $tr.Text
$idxSomeOther = $full.IndexOf("Some other ") + 1
$lenSomeOther = "Some other ".Length
$tr.Characters($idxSomeOther, $lenSomeOther).Text = ""

$full = $tr.Text
$idxOpenSource = $full.IndexOf("open-source") + 1
$tr.Characters($idxOpenSource, 1).Text = ("Some other " + $full.Substring($idxOpenSource - 1, 1))

# Expand the trailing " " after "Gurobi" into
# " (Mixed-Integer Linear Programs) ".
$full = $tr.Text
$idxGurobi = $full.IndexOf("Gurobi")
$tailStart = $idxGurobi + "Gurobi".Length + 1
$tailLen = $full.Length - $tailStart + 1
$tr.Characters($tailStart, $tailLen).Text = " (Mixed-Integer Linear Programs) "

# ---------------------------------------------------------------------
# Slide 17 notes: add a new paragraph about Gurobi's free licenses.
# ---------------------------------------------------------------------
$notesShape = $s17.NotesPage.Shapes.Item(2)
$notesTr = $notesShape.TextFrame.TextRange
$notesOrig = $notesTr.Text
$notesTr.Text = $notesOrig + "`nGurobi has free academic licenses."
